$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of this weekly block (row 368), pushing the
# existing rows 368-408 down to 370-410, then populate the two new rows with
# the new weekly price report (Primera / Segunda, Provincia de Diguillin).
$ws.Rows.Item(368).Insert()
$ws.Rows.Item(368).Insert()

# Row 368: Primera
$ws.Cells.Item(368, 1).Value = 7
$ws.Cells.Item(368, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(368, 3).Value = "Ñuble"
$ws.Cells.Item(368, 4).Value = 44918
$ws.Cells.Item(368, 5).Value = 16
$ws.Cells.Item(368, 6).Value = "Fruta"
$ws.Cells.Item(368, 7).Value = 100101
$ws.Cells.Item(368, 8).Value = "Berries"
$ws.Cells.Item(368, 9).Value = 100112025
$ws.Cells.Item(368, 10).Value = "Frutilla"
$ws.Cells.Item(368, 11).Value = "Sin especificar"
$ws.Cells.Item(368, 12).Value = "Primera"
$ws.Cells.Item(368, 13).Value = 100
$ws.Cells.Item(368, 14).Value = 6500
$ws.Cells.Item(368, 15).Value = 7000
$ws.Cells.Item(368, 16).Value = 6750
$ws.Cells.Item(368, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(368, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(368, 19).Value = 964
$ws.Cells.Item(368, 20).Value = 7

# Row 369: Segunda
$ws.Cells.Item(369, 1).Value = 7
$ws.Cells.Item(369, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(369, 3).Value = "Ñuble"
$ws.Cells.Item(369, 4).Value = 44918
$ws.Cells.Item(369, 5).Value = 16
$ws.Cells.Item(369, 6).Value = "Fruta"
$ws.Cells.Item(369, 7).Value = 100101
$ws.Cells.Item(369, 8).Value = "Berries"
$ws.Cells.Item(369, 9).Value = 100112025
$ws.Cells.Item(369, 10).Value = "Frutilla"
$ws.Cells.Item(369, 11).Value = "Sin especificar"
$ws.Cells.Item(369, 12).Value = "Segunda"
$ws.Cells.Item(369, 13).Value = 120
$ws.Cells.Item(369, 14).Value = 5000
$ws.Cells.Item(369, 15).Value = 5500
$ws.Cells.Item(369, 16).Value = 5250
$ws.Cells.Item(369, 17).Value = "`$/caja 7 kilos"
$ws.Cells.Item(369, 18).Value = "Provincia de Diguillín"
$ws.Cells.Item(369, 19).Value = 750
$ws.Cells.Item(369, 20).Value = 7
